# #5: cash & deposit done
# Update the "存款" (deposits) sheet: turn row 1 into a real header row
# (bank / deposit_type / currency / owner / total / property_category /
# category / date / legislator_name / legislator_id / source_file / index)
# and extend every data row (2-7) with the new property_category .. index
# columns. Also normalizes F6 to a real number (it was stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Row 1: header labels (B1:M1) ----
$ws.Cells.Item(1, 2).Value  = "bank"
$ws.Cells.Item(1, 3).Value  = "deposit_type"
$ws.Cells.Item(1, 4).Value  = "currency"
$ws.Cells.Item(1, 5).Value  = "owner"
$ws.Cells.Item(1, 6).Value  = "total"
$ws.Cells.Item(1, 7).Value  = "property_category"
$ws.Cells.Item(1, 8).Value  = "category"
$ws.Cells.Item(1, 9).Value  = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# ---- Data rows 2-7: add the new trailing columns G:M ----
# (column A index, property_category, category, date, legislator_name,
#  legislator_id, source_file are the same for every row)
$rows = 2..7
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value  = "deposit"
    $ws.Cells.Item($r, 8).Value  = "normal"
    $ws.Cells.Item($r, 9).Value  = "2012-03-31"
    $ws.Cells.Item($r, 10).Value = "王惠美"
    $ws.Cells.Item($r, 11).Value = 1729
    $ws.Cells.Item($r, 12).Value = "tmp99d31"
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($r, 1).Value
}

# F6 was stored as the text "4560" - normalize it to a real number like
# every other total cell in the column.
$ws.Cells.Item(6, 6).Value = 4560
